$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated case-count figures (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) for the countries whose
#     rows changed between the 13:39 and 14:56 updates on 7 Oct 2020, plus the
#     three rank swaps (Kuwait/Kazajistan, Libia/Irlanda, Santa Lucia/Nueva
#     Caledonia) that resulted from the refreshed sort order.
$updates = @(
    @{ Row=4; Cells=@{ "B"=7724725; "C"=1979; "D"=4936855; "E"=2571994; "G"=54; "H"=215876 } }
    @{ Row=20; Cells=@{ "B"=337711; "C"=468; "D"=323208; "E"=9556; "G"=24; "H"=4947 } }
    @{ Row=27; Cells=@{ "B"=278932; "C"=1906; "D"=215198; "E"=61928; "G"=9; "H"=1806 } }
    @{ Row=30; Cells=@{ "B"=149988; "C"=4989; "G"=36; "H"=6518 } }
    @{ Row=39; Cells=@{ "A"="Kuwait"; "B"=108743; "C"=475; "D"=100776; "E"=7328; "G"=7; "H"=639 } }
    @{ Row=40; Cells=@{ "A"="Kazajistan"; "B"=108362; "C"=66; "D"=103465; "E"=3151; "H"=1746 } }
    @{ Row=44; Cells=@{ "B"=101840; "C"=1046; "D"=91710; "E"=9694; "G"=1; "H"=436 } }
    @{ Row=45; Cells=@{ "B"=96677; "G"=5; "H"=5892 } }
    @{ Row=52; Cells=@{ "B"=81505; "C"=415; "D"=75683; "E"=4948; "G"=6; "H"=874 } }
    @{ Row=71; Cells=@{ "B"=41113; "C"=182; "D"=38858; "E"=1653; "G"=2; "H"=602 } }
    @{ Row=74; Cells=@{ "A"="Libia"; "B"=39513; "C"=1045; "D"=22831; "E"=16074; "G"=6; "H"=608 } }
    @{ Row=75; Cells=@{ "A"="Irlanda"; "B"=38973; "D"=23364; "E"=13798; "H"=1811 } }
    @{ Row=78; Cells=@{ "B"=30710; "C"=331; "D"=24240; "E"=5807 } }
    @{ Row=104; Cells=@{ "B"=10804; "C"=15; "E"=289; "G"=2; "H"=276 } }
    @{ Row=110; Cells=@{ "B"=9260; "C"=178; "D"=5588; "E"=3587; "G"=1; "H"=85 } }
    @{ Row=127; Cells=@{ "B"=5144; "C"=11; "D"=4885; "E"=154 } }
    @{ Row=180; Cells=@{ "B"=476; "C"=1; "E"=24 } }
    @{ Row=207; Cells=@{ "A"="Santa Lucia" } }
    @{ Row=208; Cells=@{ "A"="Nueva Caledonia" } }
)

foreach ($u in $updates) {
    foreach ($col in $u.Cells.Keys) {
        $ws.Range("$col$($u.Row)").Value = $u.Cells[$col]
    }
}

# --- Refresh the "last updated" timestamp shown on the sheet title ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Octubre de 2020 a las 14:56"

